$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (batsman) to make room for
# "ownTeam" and "oppTeam".
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data cells
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Chennai Super Kings"

$wb.Save()
